$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the current ("before") values for the columns that move, for rows 7-10,
# then write them back shifted by one row (row 10 -> row 7, row 7 -> row 8,
# row 8 -> row 9, row 9 -> row 10).

$cols = @("A", "B", "D", "E", "F", "G", "H", "Q", "R", "Z", "AB")
$rows = @(7, 8, 9, 10)

$snapshot = @{}
foreach ($r in $rows) {
    $snapshot[$r] = @{}
    foreach ($col in $cols) {
        $snapshot[$r][$col] = $ws.Range("$col$r").Value2
    }
}

# new row 7 gets old row 10's values; new row 8 gets old row 7's; etc.
$mapping = @{ 7 = 10; 8 = 7; 9 = 8; 10 = 9 }

foreach ($r in $rows) {
    $src = $mapping[$r]
    foreach ($col in $cols) {
        $ws.Range("$col$r").Value2 = $snapshot[$src][$col]
    }
}
